{"js": "// Add two runs (\"1\" with an eastAsia-hint rFonts, then \"2\" with default\n// run properties) to the end of the document's sole paragraph, matching:\n//   <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>1</w:t></w:r>\n//   <w:r><w:t>2</w:t></w:r>\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\n\n// Collapsed range at the paragraph's existing content (keeps the\n// paragraph's own attributes / <w:pPr> untouched) so the OOXML we splice\n// in only appends the two new runs instead of replacing the paragraph.\nconst insertionRange = targetParagraph.getRange(Word.RangeLocation.content);\n\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>1</w:t></w:r>' +\n  '<w:r><w:t>2</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ninsertionRange.insertOoxml(flatOpcPackage, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add two runs (\"1\" with an eastAsia-hint rFonts, then \"2\" with default\n# run properties) to the end of the document's sole paragraph, matching:\n#   <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>1</w:t></w:r>\n#   <w:r><w:t>2</w:t></w:r>\n$d = $word.ActiveDocument\n\n$targetParagraph = $d.Paragraphs(1)\n$insertionRange = $targetParagraph.Range\n\n# Shrink the range so it ends right before the paragraph mark (pilcrow);\n# otherwise InsertXML on a range that still covers the paragraph mark\n# inserts a brand-new paragraph before this one instead of appending runs\n# inside it.\n$insertionRange.MoveEnd(1, -1) | Out-Null\n\n$flatOpcPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>1</w:t></w:r>' +\n  '<w:r><w:t>2</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n$insertionRange.InsertXML($flatOpcPackage)\n"}
